$p = $ppt.ActivePresentation

# --- Remove the two slides that no longer belong in the deck ---
# (original slide6.xml "Novas Visualizacoes" and slide5.xml "Proximos Passos")
$p.Slides.Item(6).Delete()
$p.Slides.Item(5).Delete()

# --- Resize/reposition + restyle the "http://dataminas.info" slide ---
# After the two deletions it is the last slide in the deck.
$s = $p.Slides.Item($p.Slides.Count)
$shape = $s.Shapes.Item(1)

$shape.Left = 146.6519012451172
$shape.Top = 289.5
$shape.Width = 730.69633
$shape.Height = 193.87504

$tf = $shape.TextFrame
$tr = $tf.TextRange
$visible = $tr.Characters(1, 21)
$visible.Font.Size = 80
